$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sprint" column (D) data for rows 33:50 is out of order relative to the
# rest of the table (rows 6:32 are already sorted ascending by sprint). Re-sort
# the whole table body by column D ascending, as was done in Excel via
# Table > Sort, which only visibly reorders rows 33:50.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("D6:D50"), 0, 1) | Out-Null
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# Update the table's visual style.
$tbl.TableStyle = "TableStyleMedium2"

# Update the saved selection/scroll position.
$ws.Range("J13").Select()
